# Apply the "Add files via upload" edit:
#  - Rename the existing accounting rows 23-25 from category "會計" (Accounting)
#    to "預算" (Budget), keeping the rest of each row the same.
#  - Add three new rows (26-28) for category "決算" (Final Accounts / Actuals)
#    that mirror rows 23-25 but with the "actual" figures.
#  - Re-point the frozen-pane/selection back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rows 23-25: category 會計 -> 預算 -------------------------------
$ws.Range("A23").Value = "預算"
$ws.Range("A24").Value = "預算"
$ws.Range("A25").Value = "預算"

# --- 2. New rows 26-28: category 決算 -----------------------------------
$ws.Range("A26").Value = "決算"
$ws.Range("B26").Value = 113
$ws.Range("C26").Value = "工務局主管"
$ws.Range("D26").Value = "經資門合計"
$ws.Range("E26").Value = 7783220
$ws.Range("F26").Value = "113年工務局主管預算數合計8,194,228千元，執行數5,786,539千元，執行率89.70%。來源:113年工務統計年報。"

$ws.Range("A27").Value = "決算"
$ws.Range("B27").Value = 113
$ws.Range("C27").Value = "工務局主管"
$ws.Range("D27").Value = "經常門"
$ws.Range("E27").Value = 1952014
$ws.Range("F27").Value = "113年工務局主管經常門預算數2,023,521千元，執行數1,952,014千元，執行率96.47%。來源:113年工務統計年報。"

$ws.Range("A28").Value = "決算"
$ws.Range("B28").Value = 113
$ws.Range("C28").Value = "工務局主管"
$ws.Range("D28").Value = "資本門"
$ws.Range("E28").Value = 5831206
$ws.Range("F28").Value = "113年工務局主管資本門預算數6,170,707千元，執行數5,831,206千元，執行率94.50%。來源:113年工務統計年報。"

# Match row height of the sibling rows (23-25) for the newly added rows.
$ws.Rows.Item(26).RowHeight = $ws.Rows.Item(23).RowHeight
$ws.Rows.Item(27).RowHeight = $ws.Rows.Item(23).RowHeight
$ws.Rows.Item(28).RowHeight = $ws.Rows.Item(23).RowHeight

# --- 3. Scroll / selection back to the top of the frozen pane -----------
$ws.Range("A2").Select() | Out-Null
$ws.Range("D8").Select() | Out-Null
